# Update the reference for the P1/S2 and P2 monthly-cost rows: the long
# Skinner et al. (2018) citation text is replaced by its short citation
# key, "skinner2018healthcare" (ref for D3/D4 no longer duplicates the
# "None. Currently assuming same as P1/S2" footnote's neighbour string).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "skinner2018healthcare"
$ws.Range("D4").Value = "skinner2018healthcare"

# The shorter text no longer needs the extra-tall wrapped row, so let the
# rows return to their default auto-fit height.
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()

# Leave the selection where the edit finished, on D4.
$ws.Range("D4").Select()
